# feat: add 2022-Q1 data
#
# The workbook originally has two sheets: "2021-Q4" (fund holding detail)
# and "总计" (summary). This adds a new "2022-Q1" fund holding detail sheet
# (placed between the two, mirroring the "2021-Q4" sheet's layout) and
# records a new 2022-Q1 summary row at the top of "总计".

$wb = $excel.ActiveWorkbook
$detailSheet = $wb.Worksheets.Item(1)   # "2021-Q4"

# ------------------------------------------------------------------
# 1) Create the new "2022-Q1" detail sheet right after "2021-Q4" and
#    seed it from the existing sheet's layout/styles (header + the
#    A-column running index + styling all match exactly).
#
#    NOTE: sheet handles in this host resolve lazily by tab position,
#    so grab any *other* sheet reference (like "总计") only after all
#    sheet-collection mutations (Add/Move) are done.
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $detailSheet)
$newSheet.Name = "2022-Q1"

$detailSheet.Range("A1:H8").Copy($newSheet.Range("A1"))
$newSheet.Cells.Item(1, 1).ClearContents()   # A1 has no cell at all in the source layout

# ------------------------------------------------------------------
# 2) Overwrite the fund rows (B:H) with the 2022-Q1 figures. Columns
#    B,C,D,E,F,G are text values (codes/names/formatted numbers kept
#    as strings); column H is numeric. A leading apostrophe forces
#    text entry for the numeric-looking strings, then ClearFormats()
#    drops the auto-added "quote prefix" style so the cell keeps the
#    plain (unstyled) look it has in the source sheet.
# ------------------------------------------------------------------
function Set-TextCell($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

$fundRows = @(
    @("000029", "富国宏观策略灵活配置混合", "5.97", "90.27", "2.20", "0.1313", 6),
    @("010740", "汇安核心价值混合A",       "0.95", "93.68", "3.83", "0.0364", 4),
    @("005357", "富国国企改革灵活配置混合", "1.13", "87.21", "2.44", "0.0276", 5),
    @("006692", "金信消费升级股票A",       "0.62", "94.05", "3.94", "0.0244", 8),
    @("010741", "汇安核心价值混合C",       "0.36", "93.68", "3.83", "0.0138", 4),
    @("000649", "长城久鑫灵活配置混合",    "0.46", "81.08", "2.55", "0.0117", 9),
    @("006693", "金信消费升级股票C",       "0.20", "94.05", "3.94", "0.0079", 8)
)

$r = 2
foreach ($fund in $fundRows) {
    Set-TextCell $newSheet $r 2 $fund[0]
    Set-TextCell $newSheet $r 3 $fund[1]
    Set-TextCell $newSheet $r 4 $fund[2]
    Set-TextCell $newSheet $r 5 $fund[3]
    Set-TextCell $newSheet $r 6 $fund[4]
    Set-TextCell $newSheet $r 7 $fund[5]
    $newSheet.Cells.Item($r, 8).Value = $fund[6]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 3) Insert a new row at the top of "总计"'s data (row 2) for the
#    2022-Q1 summary figures, pushing the existing 2021-Q4 row down
#    to row 3. Re-apply column A's index style (copied straight from
#    the row below, which already carries it) since Insert() alone
#    pulls the header row's bold/border styling into the blank row.
#
#    Look the sheet up by name now that the sheet collection is in
#    its final shape, so the tab-position lookup lands on "总计".
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 7
$totalSheet.Cells.Item(2, 4).Value = 0.25
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Cells.Item(3, 1).Value = 1

$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)
